$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (H1:K1) to the new "Fiat"/"Asset" naming scheme.
# Order matters for where the new shared-string entries land in the table,
# so match the original authoring order: I1, J1, K1, then H1.
$ws.Range("I1").Value = "Market 1 Fiat Spot Price"
$ws.Range("J1").Value = "Market 2 Fiat Spot Price"
$ws.Range("K1").Value = "Fee Asset Fiat Spot Price"
$ws.Range("H1").Value = "Fee Asset"

# Fill in the "Fee Asset" (USD) value for every data row - row 2 already had
# it, rows 3-9 did not have an H cell at all yet.
$ws.Range("H2:H9").Value = "USD"

# Touch up the last-used selection to match the saved state.
$ws.Range("O10").Select()
